# Auto-generated edit script: updates market-price-derived columns (H-N)
# across several leve-profit sheets per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2986.25
$ws.Range("I28").Value = 2318
$ws.Range("J28").Value = 4100
$ws.Range("K28").Value = 2318
$ws.Range("L28").Value = 4100
$ws.Range("M28").Value = -1833
$ws.Range("N28").Value = -5070

$ws.Range("H88").Value = 4395.375
$ws.Range("I88").Value = 4294.5
$ws.Range("J88").Value = 4429
$ws.Range("K88").Value = 4294.5
$ws.Range("L88").Value = 4429
$ws.Range("M88").Value = -3888.5
$ws.Range("N88").Value = -5241

$ws.Range("H91").Value = 4395.375
$ws.Range("I91").Value = 4294.5
$ws.Range("J91").Value = 4429
$ws.Range("K91").Value = 4294.5
$ws.Range("L91").Value = 4429
$ws.Range("M91").Value = -2890.5
$ws.Range("N91").Value = -7237

$ws.Range("H100").Value = 6612.9546
$ws.Range("I100").Value = 2104
$ws.Range("K100").Value = 2104
$ws.Range("M100").Value = -1563

$ws.Range("H106").Value = 120320.89
$ws.Range("I106").Value = 152460.14
$ws.Range("K106").Value = 152460.14
$ws.Range("M106").Value = -151829.14

$ws.Range("H118").Value = 1352
$ws.Range("I118").Value = 1352
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 4056
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -2399
$ws.Range("N118").ClearContents()

$ws.Range("H127").Value = 1469.4286
$ws.Range("I127").Value = 857.4
$ws.Range("K127").Value = 2572.2
$ws.Range("M127").Value = 2387.8

$ws.Range("H132").Value = 20411830
$ws.Range("I132").Value = 21280304
$ws.Range("K132").Value = 63840912
$ws.Range("M132").Value = -63838382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3085.6191
$ws.Range("I2").Value = 2831.125
$ws.Range("K2").Value = 2831.125
$ws.Range("M2").Value = -2718.125

$ws.Range("H45").Value = 4057.6843
$ws.Range("I45").Value = 3071.2144
$ws.Range("J45").Value = 6819.8
$ws.Range("K45").Value = 3071.2144
$ws.Range("L45").Value = 6819.8
$ws.Range("M45").Value = -2694.2144
$ws.Range("N45").Value = -7573.8

$ws.Range("H61").Value = 7502.3105
$ws.Range("I61").Value = 4536
$ws.Range("K61").Value = 4536
$ws.Range("M61").Value = -4324

$ws.Range("H63").Value = 700400
$ws.Range("I63").Value = 800
$ws.Range("J63").Value = 1400000
$ws.Range("K63").Value = 800
$ws.Range("L63").Value = 1400000
$ws.Range("M63").Value = -114
$ws.Range("N63").Value = -1401372

$ws.Range("H66").Value = 700400
$ws.Range("I66").Value = 800
$ws.Range("J66").Value = 1400000
$ws.Range("K66").Value = 4000
$ws.Range("L66").Value = 7000000
$ws.Range("M66").Value = -568
$ws.Range("N66").Value = -7006864

$ws.Range("H116").Value = 3085.6191
$ws.Range("I116").Value = 2831.125
$ws.Range("K116").Value = 2831.125
$ws.Range("M116").Value = -537.125

$ws.Range("H136").Value = 7502.3105
$ws.Range("I136").Value = 4536
$ws.Range("K136").Value = 13608
$ws.Range("M136").Value = -11058

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3085.6191
$ws.Range("I3").Value = 2831.125
$ws.Range("K3").Value = 2831.125
$ws.Range("M3").Value = -2717.125

$ws.Range("H107").Value = 2493.4443
$ws.Range("I107").Value = 2587.7083
$ws.Range("K107").Value = 2587.7083
$ws.Range("M107").Value = -667.7082999999998

$ws.Range("H134").Value = 2258.8809
$ws.Range("I134").Value = 1724.7667
$ws.Range("J134").Value = 3594.1667
$ws.Range("K134").Value = 5174.300099999999
$ws.Range("L134").Value = 10782.5001
$ws.Range("M134").Value = -2639.300099999999
$ws.Range("N134").Value = -15852.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2386.375
$ws.Range("I16").Value = 2022.5385
$ws.Range("J16").Value = 3963
$ws.Range("K16").Value = 2022.5385
$ws.Range("L16").Value = 3963
$ws.Range("M16").Value = -1735.5385
$ws.Range("N16").Value = -4537

$ws.Range("H58").Value = 2862.4285
$ws.Range("I58").Value = 3165.4
$ws.Range("K58").Value = 3165.4
$ws.Range("M58").Value = -2962.4

$ws.Range("H113").Value = 2386.375
$ws.Range("I113").Value = 2022.5385
$ws.Range("J113").Value = 3963
$ws.Range("K113").Value = 2022.5385
$ws.Range("L113").Value = 3963
$ws.Range("M113").Value = 147.4614999999999
$ws.Range("N113").Value = -8303

$ws.Range("H132").Value = 2599.5293
$ws.Range("I132").Value = 2770.8
$ws.Range("J132").Value = 1315
$ws.Range("K132").Value = 8312.400000000001
$ws.Range("L132").Value = 3945
$ws.Range("M132").Value = -5782.400000000001
$ws.Range("N132").Value = -9005

$ws.Range("H134").Value = 2279.8
$ws.Range("I134").Value = 2133.3333
$ws.Range("K134").Value = 6399.999899999999
$ws.Range("M134").Value = -3864.999899999999

$ws.Range("H136").Value = 2862.4285
$ws.Range("I136").Value = 3165.4
$ws.Range("K136").Value = 9496.200000000001
$ws.Range("M136").Value = -6946.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 335501.34
$ws.Range("I32").Value = 2500
$ws.Range("J32").Value = 502002
$ws.Range("K32").Value = 7500
$ws.Range("L32").Value = 1506006
$ws.Range("M32").Value = -7217
$ws.Range("N32").Value = -1506572

$ws.Range("H70").Value = 2759.8
$ws.Range("I70").Value = 2633
$ws.Range("K70").Value = 7899
$ws.Range("M70").Value = -7584

$ws.Range("H73").Value = 2759.8
$ws.Range("I73").Value = 2633
$ws.Range("K73").Value = 7899
$ws.Range("M73").Value = -6807

$ws.Range("H123").Value = 1220
$ws.Range("I123").Value = 1220
$ws.Range("K123").Value = 3660
$ws.Range("M123").Value = -1210

$ws.Range("H133").Value = 7000
$ws.Range("I133").Value = 7000
$ws.Range("K133").Value = 21000
$ws.Range("M133").Value = -15940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4364.2856
$ws.Range("I113").Value = 14998
$ws.Range("J113").Value = 2592
$ws.Range("K113").Value = 14998
$ws.Range("L113").Value = 2592
$ws.Range("M113").Value = -12828
$ws.Range("N113").Value = -6932

$ws.Range("H132").Value = 3105.8374
$ws.Range("I132").Value = 3335.8064
$ws.Range("K132").Value = 10007.4192
$ws.Range("M132").Value = -7477.4192

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 20004
$ws.Range("I14").Value = 20004
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 20004
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("M14").Value = -19832

$ws.Range("H61").Value = 7903.6313
$ws.Range("I61").Value = 7591.6665
$ws.Range("J61").Value = 8438.429
$ws.Range("K61").Value = 7591.6665
$ws.Range("L61").Value = 8438.429
$ws.Range("M61").Value = -7389.6665
$ws.Range("N61").Value = -8842.429

$ws.Range("H68").Value = 2676.1365
$ws.Range("I68").Value = 2909.2727
$ws.Range("J68").Value = 2443
$ws.Range("K68").Value = 2909.2727
$ws.Range("L68").Value = 2443
$ws.Range("M68").Value = -2160.2727
$ws.Range("N68").Value = -3941

$ws.Range("H71").Value = 2676.1365
$ws.Range("I71").Value = 2909.2727
$ws.Range("J71").Value = 2443
$ws.Range("K71").Value = 14546.3635
$ws.Range("L71").Value = 12215
$ws.Range("M71").Value = -10802.3635
$ws.Range("N71").Value = -19703

$ws.Range("H82").Value = 6137.6523
$ws.Range("I82").Value = 9397.23
$ws.Range("K82").Value = 9397.23
$ws.Range("M82").Value = -9036.23

$ws.Range("H85").Value = 6137.6523
$ws.Range("I85").Value = 9397.23
$ws.Range("K85").Value = 9397.23
$ws.Range("M85").Value = -8149.23

$ws.Range("H113").Value = 7903.6313
$ws.Range("I113").Value = 7591.6665
$ws.Range("J113").Value = 8438.429
$ws.Range("K113").Value = 7591.6665
$ws.Range("L113").Value = 8438.429
$ws.Range("M113").Value = -5421.6665
$ws.Range("N113").Value = -12778.429

$ws.Range("H122").Value = 9463.777
$ws.Range("I122").Value = 8430
$ws.Range("K122").Value = 25290
$ws.Range("M122").Value = -22840

$ws.Range("H136").Value = 5169.8335
$ws.Range("I136").Value = 5225.778
$ws.Range("K136").Value = 15677.334
$ws.Range("M136").Value = -13127.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 109288
$ws.Range("I81").Value = 207577.2
$ws.Range("K81").Value = 415154.4
$ws.Range("M81").Value = -414093.4

$ws.Range("H84").Value = 109288
$ws.Range("I84").Value = 207577.2
$ws.Range("K84").Value = 2075772
$ws.Range("M84").Value = -2070468

$ws.Range("H107").Value = 890.625
$ws.Range("I107").Value = 757.875
$ws.Range("K107").Value = 2273.625
$ws.Range("M107").Value = -353.625

$ws.Range("H113").Value = 2223.3333
$ws.Range("I113").Value = 1620.4
$ws.Range("K113").Value = 4861.200000000001
$ws.Range("M113").Value = -2691.200000000001

$ws.Range("H132").Value = 1738.7778
$ws.Range("I132").Value = 1521.2858
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 4563.857400000001
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -2033.857400000001
$ws.Range("N132").Value = -12560

$ws.Range("H136").Value = 5449.9756
$ws.Range("I136").Value = 5362.303
$ws.Range("K136").Value = 16086.909
$ws.Range("M136").Value = -13536.909
